$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note (A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$cellA1 = $wsHoja1.Range("A1")
$text = $cellA1.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 15.2 = 64362.68 pesos"), "1000 Bs = 14.93 = 63208.96 pesos"
$text = $text -replace [regex]::Escape("64362.68 pesos = 15.14 = 970.17 Bs"), "63208.96 pesos = 14.87 = 952.94 Bs"
$cellA1.Value2 = $text

# --- tasas: update rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 67
$wsTasas.Range("O12").Value = 64.10299999999999
